# Updated symbol list (crypto price/volume snapshot) on Thu Feb  9 20:42:50 UTC 2023
# with GitHub Actions.
#
# For each changed cell we explicitly force a text ("@") number format before
# assigning the new value. The source cells store numeric-looking figures
# (prices) and percentages (1h volume change) as literal text strings, so
# without forcing the "@" text format Excel would silently reinterpret them
# as numbers (and percentages as fractional numbers), which would not match
# the original text content of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "312.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.37%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-7.25%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.113"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07860"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.73%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.342"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.97%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.687"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-13.01%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9241"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.69%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1075"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.55%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1783"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.58%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09038"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.76%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04412"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.38%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.207"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-16.25%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1061"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.16%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001265"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.24%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005969"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.48%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.376"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.80%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3317"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.81%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1381"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.66%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04160"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.17%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001245"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.98%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.34%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001226"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.85%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003003"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.63%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02437"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-9.23%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05317"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.32%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.008007"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.20%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1356"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.78%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007554"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.35%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.15%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008192"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.66%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3108"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.29%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006799"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.71%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.69%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003431"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.66%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004134"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "16.87%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.69%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002017"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.69%"
